# GSC export refresh: the oldest day (2025-11-16) has rolled off the
# reporting window, so drop its row from the "Chart" sheet. Excel shifts
# every subsequent row up by one, which is exactly what this commit's
# diff shows (row 2 disappears, rows 3..90 become rows 2..89, and the
# now-unused "2025-11-16" shared string is pruned automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$ws.Rows.Item(2).Delete()
